$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "platform"
$ws.Range("H2").Select() | Out-Null
